$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column C entirely (removes the search_id header and its values)
$ws.Range("C1:C3").Delete()

# Update row 2 and row 3 values
$ws.Range("A2").Value = "well_authority_number"
$ws.Range("B2").Value = "live_well_bc"

$ws.Range("A3").Value = "uwi_formatted"
$ws.Range("B3").Value = "live_well_bc"

# New rows 4-19, all referencing the live_well_bc table
$fields = @(
    "geom",
    "spud_date",
    "cumulative_oil_production_m3",
    "cumulative_gas_production_e3m3",
    "cumulative_water_production_m3",
    "cumulative_condensate_production_bbl",
    "completion_bottom_depth_m",
    "prod_ip3_oil_bbld",
    "prod_ip3_gas_mcfd",
    "full_status",
    "total_drilled_depth_m",
    "last_production_date",
    "prod_ip3_boe_boed",
    "prod_mr3_wtr_bbld",
    "prod_mr3_oil_bbld",
    "cumulative_marketable_gas_production_mcf"
)

$row = 4
foreach ($field in $fields) {
    $ws.Cells.Item($row, 1).Value = $field
    $ws.Cells.Item($row, 2).Value = "live_well_bc"
    $row = $row + 1
}
